# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells; D values
# that look like bare decimals (e.g. "589.20") get a leading apostrophe so
# Excel keeps them as literal text instead of coercing them to numbers (which
# would drop the trailing zero / precision, e.g. 589.20 -> 589.2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.697.12'
$ws.Range('E2').Value = '  +3.42%  '
$ws.Range('D3').Value = '3.689.85'
$ws.Range('E3').Value = '  +8.56%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'589.20"
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').Value = "'180.03"
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').Value = '3.679.96'
$ws.Range('E7').Value = '  +8.48%  '
$ws.Range('D8').Value = "'0.621"
$ws.Range('E8').Value = '  +4.76%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('E11').Value = '  +4.84%  '
$ws.Range('D12').Value = "'50.07"
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').Value = '4.283.22'
$ws.Range('E14').Value = '  +8.48%  '
$ws.Range('D15').Value = "'684.53"
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').Value = "'9.00"
$ws.Range('E16').Value = '  +4.66%  '
$ws.Range('D17').Value = '3.691.46'
$ws.Range('E17').Value = '  +8.52%  '
$ws.Range('D18').Value = '71.781.01'
$ws.Range('E19').Value = '  +2.30%  '
$ws.Range('D20').Value = "'18.11"
$ws.Range('E20').Value = '  +2.31%  '
$ws.Range('E21').Value = '  +3.51%  '
$ws.Range('E22').Value = '  +3.57%  '
$ws.Range('E23').Value = '  +17.64%  '
$ws.Range('D24').Value = "'17.85"
$ws.Range('E24').Value = '  +4.62%  '
$ws.Range('E25').Value = '  +2.85%  '
$ws.Range('E26').Value = '  +3.65%  '
$ws.Range('E27').Value = '  +5.39%  '
$ws.Range('D28').Value = "'10.19"
$ws.Range('E28').Value = '  +4.77%  '
$ws.Range('D29').Value = "'35.56"
$ws.Range('E29').Value = '  +6.09%  '
$ws.Range('D30').Value = "'9.25"
$ws.Range('E30').Value = '  +5.69%  '
$ws.Range('D31').Value = "'7.36"
$ws.Range('E31').Value = '  +6.81%  '
$ws.Range('D32').Value = "'4.24"
$ws.Range('E32').Value = '  +12.21%  '
$ws.Range('D33').Value = "'571.94"
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('E35').Value = '  +3.94%  '
$ws.Range('D36').Value = "'59.49"
$ws.Range('E36').Value = '  +2.47%  '
$ws.Range('D37').Value = '3.815.39'
$ws.Range('E37').Value = '  +5.75%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('E39').Value = '  +4.07%  '
$ws.Range('D40').Value = '0.0₃0782'
$ws.Range('E40').Value = '  +5.05%  '
$ws.Range('D41').Value = "'35.45"
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('D42').Value = "'3.49"
$ws.Range('E42').Value = '  +5.66%  '
$ws.Range('D43').Value = "'0.0465"
$ws.Range('E43').Value = '  +9.16%  '
$ws.Range('D44').Value = "'2.79"
$ws.Range('E44').Value = '  +3.46%  '
$ws.Range('D45').Value = "'0.354"
$ws.Range('E45').Value = '  +5.26%  '
$ws.Range('D46').Value = "'2.88"
$ws.Range('E46').Value = '  +7.68%  '
$ws.Range('D47').Value = "'3.36"
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('E48').Value = '  +4.21%  '
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('D50').Value = "'0.999"
$ws.Range('D51').Value = "'134.88"
$ws.Range('E51').Value = '  +2.83%  '
